$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("A heatmap was constructed to visualize the correlation between numerical variables, highlighting relationships like the impact of Acres on other features such as Distance to Substation. This visualization helped identify key variables that might influence the location and type of solar installations, such as proximity to substations or the size of the area.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: search text not found"
} else {
    $xml = '<w:p w14:paraId="2820FEFC" w14:textId="77777777" w:rsidR="00C811B7" w:rsidRPr="00C811B7" w:rsidRDefault="00C811B7" w:rsidP="00C811B7" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r w:rsidRPr="00C811B7"><w:lastRenderedPageBreak/><w:t>A heatmap was constructed to visualize the correlation between numerical variables, highlighting relationships like the impact of Acres on other features such as Distance to Substation. This visualization helped identify key variables that might influence the location and type of solar installations, such as proximity to substations or the size of the area.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">We can immediately see that Acres and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shape__Area</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> have a perfect correlation, suggesting that we can dismiss one of these columns in order to prevent biases for the model training. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shape__Length</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> also has a strong correlation with Acres and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shape__Area</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>th</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 0.86 which also makes sense. The distance features (to substation GTET 100, GTET 200 and CAISO) show some moderate to weak correlations with each other, with the strongest being between the distance to GTET 100 and to GTET 200 (0.70). Overall, the geographic features in the dataset like Acres and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shape__Area</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> are more strongly correlated with each other than with distance measures like GTET and CAISO. Using domain knowledge (i.e. here the &quot;business&quot; perspective) it is understandable that the solar power panel shapes and areas are strongly related to each other while the distances of the solar power installations to the different substation types indicate that installations closer to one GTET location tend to be closer to the other as well.</w:t></w:r></w:p>'
    $rng.InsertXML($xml)
    Write-Host "Paragraph replaced OK"
}
